# ECOSCOM-4453 - Create group action to export records to excel
# Changes according to the pull request comments
#
# The report template's second header cell (A2) held a literal sample
# value ("value") next to the "title" sample in A1. Per review feedback
# that sample value is removed - A2 becomes an empty, but still styled,
# cell (style index 2 is kept so the cell keeps its border/format),
# which also drops the now-unused "value" entry from the shared string
# table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").ClearContents()
